# YIELD_LV.xlsx - "Made changes to grid connection and sun profile"
#
# Updates the yearly grid-connection figures (column C) on the "Yearly"
# sheet, and tidies up the leftover border formatting on the empty
# scratch rows below the table (B8:C14), which also drops the now-blank
# edge cells (B8/C8 and B14/C14) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")

# New "sun profile" grid-connection values for 2020-2024.
$ws.Range("C2").Value = 600
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 665
$ws.Range("C5").Value = 677
$ws.Range("C6").Value = 1150

# Remove the leftover table-box border under the data (rows 8-14,
# columns B:C). The top/bottom edge cells (row 8 and row 14) become
# fully blank once their border formatting is gone, so clear them
# outright; the interior cells (rows 9-13) keep their number format.
$ws.Range("B8:C14").Borders.LineStyle = -4142
$ws.Range("B8:C8").Clear()
$ws.Range("B14:C14").Clear()

# Move the active selection to C3, matching the new cursor position.
$ws.Range("C3").Select()
